# Update "想去人数" (want-to-go count) figures for two sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1625
$wsExhibit.Range("F5").Value = 740
$wsExhibit.Range("F6").Value = 52

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1625
$wsAll.Range("F6").Value = 740
$wsAll.Range("F7").Value = 52
